$d = $word.ActiveDocument

# Paragraph 2
$p = $d.Paragraphs.Item(2)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Users (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">,  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_f_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_l_name,user_dob</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_paswd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_addr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 3
$p = $d.Paragraphs.Item(3)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Emails (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_email</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 4
$p = $d.Paragraphs.Item(4)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Telephones (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_tel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 5
$p = $d.Paragraphs.Item(5)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Profiles ( </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>,</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>profile_desciption,profile_photo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 6
$p = $d.Paragraphs.Item(6)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Photos (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>photo_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ,</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>photo_nme</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>photo_image,photo_datetime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 7
$p = $d.Paragraphs.Item(7)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Posts (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>posts_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>post_text,datetime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 8
$p = $d.Paragraphs.Item(8)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Groups (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>group_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id,group_nme,group_description</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) note: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> would be the creator’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 9
$p = $d.Paragraphs.Item(9)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>GroupMembers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>group_id,user_id</w:t></w:r><w:r><w:t>,member_status</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) note: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>member_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> either being creator, content creator or member</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 10
$p = $d.Paragraphs.Item(10)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>GroupPosts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>group_id,post_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 11
$p = $d.Paragraphs.Item(11)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Friends ( </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>user_id,</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>friend_id</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>,</w:t></w:r><w:r><w:t>friend_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 12
$p = $d.Paragraphs.Item(12)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Comment(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>comment_id,</w:t></w:r><w:r><w:t>user_id,post_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>,</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>comment_text,comment_datetime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)

# Paragraph 13
$p = $d.Paragraphs.Item(13)
$start = $p.Range.Start
$end = $p.Range.End
$sub = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Guest (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>guest_ID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>guest_nme</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>guest_email</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>guest_addr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gues</w:t></w:r><w:r><w:t>z</w:t></w:r><w:r><w:t>t_paswd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>guest_tel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sub.InsertXML($xml)
